$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (report week / volume number)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/30/2024  Through  10/6/2024"

# ---------------------------------------------------------------------------
# Cells that flip between a numeric value and the literal text placeholders
# ("0" / "***.*") need their donor cell's exact style + type copied over so
# the shared-string / numeric typing matches exactly. Donor cells below are
# never themselves touched by this edit, so they stay stable reference
# points for the whole script.
#   C14 -> text "0"     (style 14)
#   E14 -> text "***.*" (style 14)
#   I14 -> number 1     (style 15)
#   G31 -> number 2     (style 15)
#   M14 -> number -50   (style 16)
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("I14").Copy($ws.Range("C18"))
$ws.Range("G31").Copy($ws.Range("D18"))
$ws.Range("M14").Copy($ws.Range("E18"))
$ws.Range("G31").Copy($ws.Range("D22"))
$ws.Range("M14").Copy($ws.Range("E22"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("C31"))
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$ws.Range("L15").Value = -60

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -18.181818181818
$ws.Range("I16").Value = 100
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = -11.504424778761
$ws.Range("L16").Value = 23.456790123456
$ws.Range("M16").Value = -20.634920634920
$ws.Range("N16").Value = -83.193277310924

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -17.647058823529
$ws.Range("I17").Value = 138
$ws.Range("J17").Value = 142
$ws.Range("K17").Value = -2.816901408450
$ws.Range("L17").Value = 14.049586776859
$ws.Range("M17").Value = 53.333333333333
$ws.Range("N17").Value = -37.837837837837

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 42
$ws.Range("J18").Value = 75
$ws.Range("K18").Value = -44
$ws.Range("L18").Value = -57.142857142857
$ws.Range("M18").Value = -28.813559322033
$ws.Range("N18").Value = -91.025641025641

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 62.5
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -3.125
$ws.Range("I19").Value = 266
$ws.Range("J19").Value = 308
$ws.Range("K19").Value = -13.636363636363
$ws.Range("L19").Value = -0.374531835205
$ws.Range("M19").Value = 30.392156862745
$ws.Range("N19").Value = -48.449612403100

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 31
$ws.Range("J20").Value = 83
$ws.Range("K20").Value = -62.650602409638
$ws.Range("L20").Value = -44.642857142857
$ws.Range("M20").Value = 55
$ws.Range("N20").Value = -89.419795221843

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 26.315789473684
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = -13.698630136986
$ws.Range("I21").Value = 584
$ws.Range("J21").Value = 728
$ws.Range("K21").Value = -19.780219780219
$ws.Range("L21").Value = -8.607198748043
$ws.Range("M21").Value = 13.840155945419
$ws.Range("N21").Value = -72.517647058823

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 17
$ws.Range("J22").Value = 29
$ws.Range("K22").Value = -41.379310344827
$ws.Range("L22").Value = 13.333333333333
$ws.Range("M22").Value = -26.086956521739

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 4
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 30
$ws.Range("I23").Value = 105
$ws.Range("J23").Value = 118
$ws.Range("K23").Value = -11.016949152542
$ws.Range("L23").Value = 15.384615384615
$ws.Range("M23").Value = 69.354838709677

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 6
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 31
$ws.Range("G24").Value = 45
$ws.Range("H24").Value = -31.111111111111
$ws.Range("I24").Value = 312
$ws.Range("J24").Value = 394
$ws.Range("K24").Value = -20.812182741116
$ws.Range("L24").Value = -21.410579345088
$ws.Range("M24").Value = -35.802469135802

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = -33.333333333333
$ws.Range("I25").Value = 54
$ws.Range("J25").Value = 134
$ws.Range("K25").Value = -59.701492537313
$ws.Range("L25").Value = -60.294117647058

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 53.846153846153
$ws.Range("I26").Value = 185
$ws.Range("J26").Value = 183
$ws.Range("K26").Value = 1.092896174863
$ws.Range("L26").Value = 12.121212121212
$ws.Range("M26").Value = -18.502202643171

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = 10
$ws.Range("L27").Value = -50

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -44.444444444444
$ws.Range("J28").Value = 38
$ws.Range("K28").Value = -18.421052631578
$ws.Range("L28").Value = -18.421052631578

Write-Output "edits applied"
